$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..N (column G is untouched, remains 0), rows 2..25.
$updates = @(
    @("B2", 0.956140796269068),
    @("C2", 0.15900708755143),
    @("D2", 0.005043315098229417),
    @("E2", 0.01340457443459364),
    @("F2", 3.154037444716565),
    @("H2", 0.04311626980804162),
    @("I2", 0.05721842803638477),
    @("J2", 0.0267246236506935),
    @("K2", 2.046771444283237),
    @("L2", 0.05988879959687665),
    @("M2", 1.053670662352886),
    @("N2", 0.1775815057853691),
    @("B3", 0.8366025064731559),
    @("C3", 0.1381913326310524),
    @("D3", 0.005109710874465323),
    @("E3", 0.01171277533676474),
    @("F3", 2.930293286779076),
    @("H3", 0.05097056615325268),
    @("I3", 0.0665632146195887),
    @("J3", 0.0267246236506935),
    @("K3", 1.912285751662537),
    @("L3", 0.0534943589267769),
    @("M3", 0.9212702128117769),
    @("N3", 0.1566961188261189),
    @("B4", 0.763420788237795),
    @("C4", 0.1257434158767552),
    @("D4", 0.005138051146845068),
    @("E4", 0.01068859273770784),
    @("F4", 2.792187522793114),
    @("H4", 0.05629584080980887),
    @("I4", 0.07289198219109938),
    @("J4", 0.0267246236506935),
    @("K4", 1.828835017825043),
    @("L4", 0.04954824554659965),
    @("M4", 0.8405608898876835),
    @("N4", 0.143971709379791),
    @("B5", 0.7335000695912584),
    @("C5", 0.1209681150842385),
    @("D5", 0.005131237174371694),
    @("E5", 0.0102775579787453),
    @("F5", 2.731236135821334),
    @("H5", 0.05860875302573554),
    @("I5", 0.07574575787719651),
    @("J5", 0.0267246236506935),
    @("K5", 1.791186946846466),
    @("L5", 0.0478389961136898),
    @("M5", 0.8080058530143788),
    @("N5", 0.1389749554129907),
    @("B6", 0.7283594327541039),
    @("C6", 0.1204498192676056),
    @("D6", 0.00511190855197563),
    @("E6", 0.0102135088497004),
    @("F6", 2.715635966266476),
    @("H6", 0.0590253023051277),
    @("I6", 0.076395289999863),
    @("J6", 0.0267246236506935),
    @("K6", 1.78075133048398),
    @("L6", 0.04743904820827893),
    @("M6", 0.8028570305472158),
    @("N6", 0.1383503141173179),
    @("B7", 0.7625345596210877),
    @("C7", 0.1264173385456928),
    @("D7", 0.005088063135610721),
    @("E7", 0.01069404670007046),
    @("F7", 2.776417088640912),
    @("H7", 0.05639302142459002),
    @("I7", 0.07336460328171768),
    @("J7", 0.0267246236506935),
    @("K7", 1.816921564504042),
    @("L7", 0.04920764485095219),
    @("M7", 0.8407997686515785),
    @("N7", 0.144461931755103),
    @("B8", 0.9142102448249148),
    @("C8", 0.1527385124681757),
    @("D8", 0.00500297989524956),
    @("E8", 0.0128323031517672),
    @("F8", 3.057200299437056),
    @("H8", 0.04579435085154282),
    @("I8", 0.06083217756445336),
    @("J8", 0.0267246236506935),
    @("K8", 1.985442811865653),
    @("L8", 0.05726083636836776),
    @("M8", 1.008771703620283),
    @("N8", 0.1711016927548386),
    @("B9", 1.214761326732173),
    @("C9", 0.2059951658961126),
    @("D9", 0.004825187757835803),
    @("E9", 0.01715205808953435),
    @("F9", 3.629702803030114),
    @("H9", 0.0291028302534202),
    @("I9", 0.0405163454041757),
    @("J9", 0.0267246236506935),
    @("K9", 2.329807203136696),
    @("L9", 0.07354401980651204),
    @("M9", 1.342411163224881),
    @("N9", 0.2232530569038573),
    @("B10", 1.437320151505389),
    @("C10", 0.2476791300692298),
    @("D10", 0.004457285504360087),
    @("E10", 0.01949380530071743),
    @("F10", 3.964805645037075),
    @("H10", 0.0203144793315313),
    @("I10", 0.02939875534950254),
    @("J10", 0.0267246236506935),
    @("K10", 2.521751134168483),
    @("L10", 0.08141781694210337),
    @("M10", 1.592080434545409),
    @("N10", 0.2545962256738648),
    @("B11", 1.537237154997058),
    @("C11", 0.2674376553424338),
    @("D11", 0.00336896022937383),
    @("E11", 0.01359581140170096),
    @("F11", 3.457402381428949),
    @("H11", 0.03844125354471473),
    @("I11", 0.02786860026064542),
    @("J11", 0.0267246236506935),
    @("K11", 2.161364983301951),
    @("L11", 0.05651790406298218),
    @("M11", 1.705205817691223),
    @("N11", 0.2042132700840256),
    @("B12", 1.575184204089481),
    @("C12", 0.2738897386972212),
    @("D12", 0.003083239896706225),
    @("E12", 0.01003935200219841),
    @("F12", 3.013041851711122),
    @("H12", 0.07774722060106853),
    @("I12", 0.02789529532480994),
    @("J12", 0.0267246236506935),
    @("K12", 1.860372498987886),
    @("L12", 0.04294612228613381),
    @("M12", 1.746462363547749),
    @("N12", 0.1594931861536324),
    @("B13", 1.566451161026777),
    @("C13", 0.2714899921468827),
    @("D13", 0.003274135213283458),
    @("E13", 0.008084279965379748),
    @("F13", 2.56944121461791),
    @("H13", 0.135152698745685),
    @("I13", 0.02965752003883804),
    @("J13", 0.02672462365057982),
    @("K13", 1.57063167185548),
    @("L13", 0.03678467280069153),
    @("M13", 1.735367948599702),
    @("N13", 0.1177376203333438),
    @("B14", 1.539106317133303),
    @("C14", 0.2660677086321073),
    @("D14", 0.003676294897365651),
    @("E14", 0.007705019766997001),
    @("F14", 2.262510898395192),
    @("H14", 0.1861310804173115),
    @("I14", 0.03183932013758373),
    @("J14", 0.02672462365057982),
    @("K14", 1.375353075029423),
    @("L14", 0.03697423085458773),
    @("M14", 1.704214975412299),
    @("N14", 0.09097433416223311),
    @("B15", 1.522307315369432),
    @("C15", 0.2631963866126341),
    @("D15", 0.003843105580115136),
    @("E15", 0.00770077950510295),
    @("F15", 2.180517466958747),
    @("H15", 0.1993861273887489),
    @("I15", 0.03296093852107784),
    @("J15", 0.02672462365057982),
    @("K15", 1.324383619221514),
    @("L15", 0.03749446254114197),
    @("M15", 1.685798060150944),
    @("N15", 0.08459185689875426),
    @("B16", 1.427717374508489),
    @("C16", 0.2458727567494918),
    @("D16", 0.003905940448971634),
    @("E16", 0.007157295094823768),
    @("F16", 2.12884978632114),
    @("H16", 0.1903065485740711),
    @("I16", 0.03750546046356185),
    @("J16", 0.02672462365057982),
    @("K16", 1.301682626030313),
    @("L16", 0.03564643560948788),
    @("M16", 1.580515164348327),
    @("N16", 0.08122936943355796),
    @("B17", 1.369859442481442),
    @("C17", 0.2355954096266686),
    @("D17", 0.003711538602633624),
    @("E17", 0.006760190342802452),
    @("F17", 2.247096441320124),
    @("H17", 0.1545121134063265),
    @("I17", 0.03990310692482968),
    @("J17", 0.02672462365057982),
    @("K17", 1.384024150464143),
    @("L17", 0.03326719592992955),
    @("M17", 1.516593918822082),
    @("N17", 0.09218279335428292),
    @("B18", 1.337111613067151),
    @("C18", 0.2293720138151372),
    @("D18", 0.003437278992257653),
    @("E18", 0.007381777294487968),
    @("F18", 2.551466892361788),
    @("H18", 0.101931622860981),
    @("I18", 0.04003802812270507),
    @("J18", 0.02672462365057982),
    @("K18", 1.586759257969135),
    @("L18", 0.03437282184863388),
    @("M18", 1.479892550895983),
    @("N18", 0.1191743897711035),
    @("B19", 1.325764883468111),
    @("C19", 0.2282458402223995),
    @("D19", 0.003535427801430568),
    @("E19", 0.01019836242551975),
    @("F19", 2.989571040433759),
    @("H19", 0.05506183589662328),
    @("I19", 0.0389062980604864),
    @("J19", 0.02672462365057982),
    @("K19", 1.877653834019341),
    @("L19", 0.04435627284351007),
    @("M19", 1.468823316359419),
    @("N19", 0.1622774651128438),
    @("B20", 1.377013690421165),
    @("C20", 0.238903312044215),
    @("D20", 0.004425165629206917),
    @("E20", 0.01885925517707854),
    @("F20", 3.827217467502976),
    @("H20", 0.02255227809694871),
    @("I20", 0.0332341846476174),
    @("J20", 0.02672462365057982),
    @("K20", 2.43396551948193),
    @("L20", 0.07816813284429003),
    @("M20", 1.528173249544636),
    @("N20", 0.2477553658803231),
    @("B21", 1.549138581477763),
    @("C21", 0.2714378978198511),
    @("D21", 0.004380755906254308),
    @("E21", 0.02220924376697653),
    @("F21", 4.211778126924912),
    @("H21", 0.01583050657915441),
    @("I21", 0.02511597596692283),
    @("J21", 0.02672462365057982),
    @("K21", 2.667934989884159),
    @("L21", 0.09058470969535648),
    @("M21", 1.72149065560393),
    @("N21", 0.2841210390492108),
    @("B22", 1.662973320778264),
    @("C22", 0.2922536343155286),
    @("D22", 0.004322228143974316),
    @("E22", 0.02392950947240458),
    @("F22", 4.444147681205152),
    @("H22", 0.01232351856274855),
    @("I22", 0.02025754594460061),
    @("J22", 0.02672462365057982),
    @("K22", 2.809078938243715),
    @("L22", 0.09717737031336071),
    @("M22", 1.848125867237684),
    @("N22", 0.3032588286082216),
    @("B23", 1.60281643647258),
    @("C23", 0.2801973235767434),
    @("D23", 0.004403956139819698),
    @("E23", 0.02299499646113112),
    @("F23", 4.337874080541013),
    @("H23", 0.01410047062423891),
    @("I23", 0.02237031077245),
    @("J23", 0.02672462365057982),
    @("K23", 2.747332913524517),
    @("L23", 0.09405152882647627),
    @("M23", 1.779677073893964),
    @("N23", 0.2923444751149873),
    @("B24", 1.37526336458518),
    @("C24", 0.2370969368251394),
    @("D24", 0.004612496626362761),
    @("E24", 0.01955618956714034),
    @("F24", 3.90816204176619),
    @("H24", 0.0221305744377116),
    @("I24", 0.03244585721625537),
    @("J24", 0.02672462365057982),
    @("K24", 2.491669442223468),
    @("L24", 0.08168486398587049),
    @("M24", 1.524071910539249),
    @("N24", 0.2524724671905005),
    @("B25", 1.132192372193288),
    @("C25", 0.1925738188881212),
    @("D25", 0.004800638877207497),
    @("E25", 0.01598289366406114),
    @("F25", 3.44809866848226),
    @("H25", 0.03322520670206552),
    @("I25", 0.04612562587991764),
    @("J25", 0.02672462365057982),
    @("K25", 2.216596152124666),
    @("L25", 0.06854579678202199),
    @("M25", 1.252633865324725),
    @("N25", 0.2100548682253844)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
